# Finished Extensor 10 mm data collection
# Fill in the remaining test-run columns (E:K) for the Extensor 10mm
# results table, and move the active selection to K14 to reflect where
# data entry stopped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtTest10mm")

# Row 6 - Load (N)
$ws.Range("E6").Value = 10.593
$ws.Range("F6").Value = 11.395
$ws.Range("G6").Value = 10.655
$ws.Range("H6").Value = 7.7816
$ws.Range("I6").Value = 3.958
$ws.Range("J6").Value = 1.909
$ws.Range("K6").Value = 1.22

# Row 7 - Knee angle
$ws.Range("E7").Value = 93
$ws.Range("F7").Value = 70
$ws.Range("G7").Value = 53.5
$ws.Range("H7").Value = 49
$ws.Range("I7").Value = 21
$ws.Range("J7").Value = 9
$ws.Range("K7").Value = 2

# Row 8 - MA
$ws.Range("E8").Value = 36.1
$ws.Range("F8").Value = 39
$ws.Range("G8").Value = 34
$ws.Range("H8").Value = 38.5
$ws.Range("I8").Value = 36.4
$ws.Range("J8").Value = 33.3
$ws.Range("K8").Value = 34.6

# Row 9 - ICR (spline)
$ws.Range("E9").Value = 34
$ws.Range("F9").Value = 42.5
$ws.Range("G9").Value = 40
$ws.Range("H9").Value = 41.5
$ws.Range("I9").Value = 38
$ws.Range("J9").Value = 38
$ws.Range("K9").Value = 35

# Row 10 - Tibia origin (matlab)
$ws.Range("E10").Value = 510
$ws.Range("F10").Value = 495
$ws.Range("G10").Value = 485
$ws.Range("H10").Value = 475
$ws.Range("I10").Value = 460
$ws.Range("J10").Value = 450
$ws.Range("K10").Value = 450

# Row 13 - Load cell angle (tibia)
$ws.Range("E13").Value = 40
$ws.Range("F13").Value = 40
$ws.Range("G13").Value = 45
$ws.Range("H13").Value = 45
$ws.Range("I13").Value = 40
$ws.Range("J13").Value = 40
$ws.Range("K13").Value = 41

# Row 15 (Torque) recalculates automatically from the shared formula.

# Reflect where data entry left off.
$ws.Range("K14").Select() | Out-Null
